# Update column G ("K") values for rows 2-6 as part of regenerating
# save_data to use K instead of Strike# (recalculated std/mean, s_vals).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 4
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 1
